# Characters.xlsx edit script
# ---------------------------------------------------------------
# Reproduces the "torn categories" data fix: Bub's row (previously the
# last entry of the Level-4 bucket, row 86) is moved up to become the
# first entry of the Level-1 bucket (row 57). Every row between the old
# and new position (57-85) shifts down by one row, which is the
# observable effect of a classic "cut row 86, insert before row 57"
# spreadsheet edit.
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 57
$lastRow  = 86

# Snapshot the current (pre-edit) values for columns A (Character),
# B (Level), C (Class/Subclass) and E (Game) across the affected rows.
$origA = @{}
$origB = @{}
$origC = @{}
$origE = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $origA[$r] = $ws.Cells.Item($r, 1).Value2
    $origB[$r] = $ws.Cells.Item($r, 2).Value2
    $origC[$r] = $ws.Cells.Item($r, 3).Value2
    $origE[$r] = $ws.Cells.Item($r, 5).Value2
}

# Rows 58..86 take on the A/B/C/E values that used to live one row above
# them (the classic "shift down" that happens when a row is inserted
# above and the old last row is removed from the bottom).
for ($r = $lastRow; $r -ge ($firstRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 1).Value = $origA[$src]
    $ws.Cells.Item($r, 2).Value = $origB[$src]
    $ws.Cells.Item($r, 3).Value = $origC[$src]
    $ws.Cells.Item($r, 5).Value = $origE[$src]
}

# Row 57 becomes the moved row (old row 86 / Bub): Character, Class and
# Game move with it, but the Level (column B) keeps the destination
# row's original value, i.e. Bub is now unlocked at Level 1 instead of
# Level 4.
$ws.Cells.Item($firstRow, 1).Value = $origA[$lastRow]
$ws.Cells.Item($firstRow, 2).Value = $origB[$firstRow]
$ws.Cells.Item($firstRow, 3).Value = $origC[$lastRow]
$ws.Cells.Item($firstRow, 5).Value = $origE[$lastRow]

# ---------------------------------------------------------------
# Sheet view / selection housekeeping to match the saved workbook state
# ---------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 64
$win.ScrollColumn = 1
$ws.Range("C75").Select()
